$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New columns J (10) and K (11): widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 15.6666666667
$ws.Columns.Item(11).ColumnWidth = 14.6666666667

# ---------------------------------------------------------------------------
# 2. Header area J2:K2 ("Speed Increase", merged) and J3/K3 ("Main Window, %"
#    / "Input Window, %"), re-using formatting from the existing matching
#    header cells so the same style indexes get reused.
# ---------------------------------------------------------------------------
$ws.Range("G2:H2").Copy()
$ws.Range("J2:K2").PasteSpecial(-4122)
$ws.Range("J2").Value = "Speed Increase"

$ws.Range("G3:H3").Copy()
$ws.Range("J3:K3").PasteSpecial(-4122)
$ws.Range("J3").Value = "Main Window, %"
$ws.Range("K3").Value = "Input Window, %"

$ws.Range("J2:K2").Merge()

# ---------------------------------------------------------------------------
# 3. Data area J4:K103 - copy the bordered number-cell style (from B4) across
#    the whole block in one shot, then apply the percentage number format.
# ---------------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("J4:K105").PasteSpecial(-4122)
$ws.Range("J4:K105").NumberFormat = "0%"

# Row 4 gets its own (non shared) formulas.
$ws.Range("J4").Formula = "=(C4-G4)/C4"
$ws.Range("K4").Formula = "=(D4-H4)/D4"

# Rows 5-68 share one formula group, rows 69-103 share another - matching
# how Excel itself groups a single fill operation into one shared formula.
$ws.Range("J5:K68").Formula = "=(C5-G5)/C5"
$ws.Range("J69:K103").Formula = "=(C69-G69)/C69"

# Row 105 - per-column averages (kept as individual, non shared formulas).
$ws.Range("J105").Formula = "=AVERAGE(J4:J103)"
$ws.Range("K105").Formula = "=AVERAGE(K4:K103)"

# ---------------------------------------------------------------------------
# 4. Reposition the two charts that overlapped the new columns.
# ---------------------------------------------------------------------------
$co1 = $ws.ChartObjects(1)
$co1.Left = 879.8646681225393
$co1.Top = 24.496614173228345
$co1.Width = 1375.2587204724412
$co1.Height = 530.4000000000001

$co2 = $ws.ChartObjects(2)
$co2.Left = 881.1130933193898
$co2.Top = 570.402283464567
$co2.Width = 1375.258720472441
$co2.Height = 530.4920472440946

# ---------------------------------------------------------------------------
# 5. View: zoom back to 100%, scroll so column G is left-most, select K14.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("K14").Select()

Write-Output "done"
